# Update the "Point of Contact" block of the package doc template:
#  - the contact name becomes the {{point_of_contact}} merge placeholder
#  - the old title / email / blank line that followed it is removed entirely

$d = $word.ActiveDocument

# 1. Replace the literal contact name with the templated placeholder.
$d.Content.Find.Execute("Jennie Rice", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{{point_of_contact}}", 2)

# 2. Drop the paragraph that used to carry the title, line break and email
#    address — it's gone completely now that the name is a placeholder.
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*IM3 Principal Investigator*jennie.rice@pnnl.gov*") {
        $p.Range.Delete()
        break
    }
}
